# Apply cell updates from diff (rows 2-51 of cryptos sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.633.60"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "2.508.34"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'574.17"
$ws.Range("E5").Value = "  -0.84%  "
$ws.Range("D6").Value = "'166.95"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "'0.515"
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("D9").Value = "2.507.78"
$ws.Range("E9").Value = "  +0.34%  "
$ws.Range("D10").Value = "'0.161"
$ws.Range("E10").Value = "  +1.78%  "
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("E12").Value = "  +6.36%  "
$ws.Range("E13").Value = "  +1.83%  "
$ws.Range("D14").Value = "2.971.31"
$ws.Range("E14").Value = "  +0.33%  "
$ws.Range("D15").Value = "69.395.17"
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("E16").Value = "  +1.47%  "
$ws.Range("D17").Value = "'24.84"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").Value = "2.519.22"
$ws.Range("E18").Value = "  +0.70%  "
$ws.Range("D19").Value = "'11.30"
$ws.Range("E19").Value = "  -0.67%  "
$ws.Range("D20").Value = "'7.59"
$ws.Range("E20").Value = "  -1.66%  "
$ws.Range("D21").Value = "'351.70"
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").Value = "'3.91"
$ws.Range("E22").Value = "  -0.77%  "
$ws.Range("D23").Value = "'1.96"
$ws.Range("E23").Value = "  -0.66%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "'0.999"
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'70.97"
$ws.Range("E25").Value = "  +2.63%  "
$ws.Range("D26").Value = "'3.95"
$ws.Range("E26").Value = "  -1.72%  "
$ws.Range("D27").Value = "'8.84"
$ws.Range("E27").Value = "  -2.22%  "
$ws.Range("D28").Value = "2.666.00"
$ws.Range("D29").Value = "'0.997"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "0.0₃0890"
$ws.Range("E30").Value = "  -1.31%  "
$ws.Range("D31").Value = "'7.88"
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("D32").Value = "'461.25"
$ws.Range("E32").Value = "  -4.52%  "
$ws.Range("D33").Value = "'1.22"
$ws.Range("E33").Value = "  -5.50%  "
$ws.Range("E34").Value = "  -1.16%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "'160.05"
$ws.Range("E36").Value = "  +5.22%  "
$ws.Range("E37").Value = "  +1.39%  "
$ws.Range("E38").Value = "  +1.05%  "
$ws.Range("D39").Value = "'18.50"
$ws.Range("E39").Value = "  -0.29%  "
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").Value = "'0.318"
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("E42").Value = "  -1.81%  "
$ws.Range("E43").Value = "  -1.59%  "
$ws.Range("D44").Value = "'38.13"
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("E45").Value = "  -5.09%  "
$ws.Range("D46").Value = "'1.09"
$ws.Range("E46").Value = "  -6.70%  "
$ws.Range("D47").Value = "'142.33"
$ws.Range("E47").Value = "  -0.72%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "'0.521"
$ws.Range("E48").Value = "  -1.55%  "
$ws.Range("B49").Value = "Filecoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D49").Value = "'3.47"
$ws.Range("E49").Value = "  -1.72%  "
$ws.Range("E50").Value = "  +0.63%  "
$ws.Range("D51").Value = "'5.78"
$ws.Range("E51").Value = "  +3.01%  "
